$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignIn")
$ws.Activate()

# Fill A3 with the same Test case name value as A2 ("Verify_ebfs_Login")
$ws.Range("A3").Value = $ws.Range("A2").Value2

# Move the active selection to A3, matching the saved cursor position
$ws.Range("A3").Select()
